$d = $word.ActiveDocument
$d.Content.Find.Execute("using  what", $true, $false, $false, $false, $false, $true, 1, $false, "using what", 2)
